$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "article 82 is live": the rotating blog-post slots (I7 -> ser 79, E7 -> ser 80,
# C7 -> ser 81) each roll forward by one ("ser 79" is retired and its shared
# string is dropped, which bumps the later "ser 80"/"ser 81" strings down an
# index and implicitly makes I7/E7 display what used to be in E7/C7); C7 then
# shows the brand-new "ser 82" post. D7 (this week's meetup placeholder) is
# replaced with the freshly announced meetup.
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 80"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 81"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 82"

$newMeetup = "type: meetup`nwidth: 2`nheight: 1`nh3: Meetup coming in`ndate: 2020,3,5,10,30,0,0`nbutton.default: Speak*goto(`"https://forms.gle/dyydXFRSsKzeH4hZ6`")`nbutton.default: Attend*goto(`"https://youtu.be/vscn-HP932E`")`nbutton.default: Details*goto(`"https://www.meetup.com/techshek/events/269581504/`")"
$ws.Range("D7").Value = $newMeetup

# Scroll/selection moved up-and-left by one cell (B7->B6, I7->E7) in the saved view.
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 2
$ws.Range("E7").Select() | Out-Null
